$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price-table refresh: columns D (Price) and E (Volume/1h%)
# are stored as literal text in the sheet (not numbers), so trailing zeros
# ("99.20") and thousands-dot-grouped prices ("45.747.69") survive exactly.
#
# Plain "$ws.Range(ref).Value = '123.45'" would get auto-coerced by Excel
# into the Number type for anything that parses as a plain decimal, losing
# formatting (trailing zeros) and changing the stored cell type. To keep an
# un-ambiguous literal string (and leave the cell's existing/default style
# untouched), route those values through a scratch cell: put a formula that
# evaluates to the literal text, copy it, Paste-Special into the destination
# (this carries the Text type without touching NumberFormat/style), then
# clear the scratch cell again.
function Set-LiteralText($ref, $value) {
    $scratch = $ws.Range("ZZ1000")
    $scratch.Formula = '="' + $value + '"'
    $scratch.Copy()
    $ws.Range($ref).PasteSpecial()
    $scratch.ClearContents()
}

$ws.Range("D2").Value = '45.747.69'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '2.373.24'
$ws.Range("E3").Value = '  +3.00%  '
$ws.Range("E4").Value = '  -0.11%  '
Set-LiteralText "D5" '299.36'
$ws.Range("E5").Value = '  -2.09%  '
Set-LiteralText "D6" '97.29'
$ws.Range("E6").Value = '  -4.69%  '
$ws.Range("E7").Value = '  -1.91%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -4.54%  '
Set-LiteralText "D10" '33.85'
$ws.Range("E10").Value = '  -8.97%  '
$ws.Range("E11").Value = '  -1.86%  '
Set-LiteralText "D12" '7.04'
$ws.Range("E12").Value = '  -5.99%  '
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").Value = '2.736.65'
$ws.Range("D15").Value = '2.368.77'
$ws.Range("E15").Value = '  +2.80%  '
$ws.Range("E16").Value = '  -0.51%  '
Set-LiteralText "D17" '13.63'
$ws.Range("E17").Value = '  -3.05%  '
$ws.Range("D18").Value = '45.681.55'
$ws.Range("E18").Value = '  -2.46%  '
Set-LiteralText "D19" '12.66'
$ws.Range("E19").Value = '  -6.31%  '
$ws.Range("E20").Value = '  -1.15%  '
Set-LiteralText "D21" '6.01'
$ws.Range("E21").Value = '  -1.32%  '
Set-LiteralText "D22" '66.70'
$ws.Range("E22").Value = '  -0.33%  '
Set-LiteralText "D23" '242.53'
$ws.Range("E23").Value = '  -3.23%  '
$ws.Range("E24").Value = '  -6.35%  '
$ws.Range("E25").Value = '  -0.03%  '
Set-LiteralText "D26" '1.90'
$ws.Range("E26").Value = '  -3.37%  '
Set-LiteralText "D27" '38.35'
$ws.Range("E27").Value = '  -13.47%  '
$ws.Range("E28").Value = '  -3.54%  '
Set-LiteralText "D29" '9.63'
$ws.Range("E29").Value = '  -3.52%  '
Set-LiteralText "D30" '3.76'
$ws.Range("E30").Value = '  +16.26%  '
Set-LiteralText "D31" '20.93'
$ws.Range("E31").Value = '  +3.41%  '
Set-LiteralText "D32" '2.76'
$ws.Range("E32").Value = '  -4.31%  '
Set-LiteralText "D33" '5.48'
$ws.Range("E33").Value = '  -5.09%  '
Set-LiteralText "D34" '146.22'
Set-LiteralText "D35" '0.0766'
$ws.Range("E35").Value = '  -4.77%  '
Set-LiteralText "D36" '0.111'
$ws.Range("E36").Value = '  -2.93%  '
Set-LiteralText "D37" '1.91'
$ws.Range("E37").Value = '  +5.60%  '
Set-LiteralText "D38" '0.116'
$ws.Range("E38").Value = '  -3.04%  '
Set-LiteralText "D39" '15.11'
$ws.Range("E39").Value = '  -7.15%  '
$ws.Range("E40").Value = '  -7.48%  '
Set-LiteralText "D41" '0.0295'
$ws.Range("E41").Value = '  -3.46%  '
Set-LiteralText "D42" '3.19'
$ws.Range("E42").Value = '  -8.20%  '
$ws.Range("D43").Value = '1.945.73'
$ws.Range("E43").Value = '  +4.71%  '
$ws.Range("E44").Value = '  -0.05%  '
Set-LiteralText "D45" '93.71'
$ws.Range("E45").Value = '  +5.83%  '
$ws.Range("E46").Value = '  -10.77%  '
Set-LiteralText "D47" '8.47'
$ws.Range("E47").Value = '  +6.36%  '
Set-LiteralText "D48" '99.20'
$ws.Range("E48").Value = '  +2.32%  '
$ws.Range("E49").Value = '  -7.71%  '
$ws.Range("D50").Value = '2.605.97'
$ws.Range("E50").Value = '  +2.92%  '
Set-LiteralText "D51" '68.34'
$ws.Range("E51").Value = '  -8.97%  '
